# Update Supply values on the active sheet (Warehouse_FEMA_4 / Supply_info).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1500
$ws.Range("B3").Value = 1500
$ws.Range("B4").Value = 1500

# Move / persist the active selection to G5, as in the saved workbook view.
$ws.Range("G5").Select()
